$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 21 (Anna) and row 22 (Stephan) for LeetCode #88 "Merge Sorted Array",
# then add two brand-new rows (23/24) for #176 "Second Highest Salary", and a
# trailing blank formatted row (25) -- mirroring the existing table's rhythm of
# one row per person per problem.

$ws.Range("D21").Value = "88. Merge Sorted Array"
$ws.Range("D22").Value = "88. Merge Sorted Array"
$ws.Range("E22").Value = "2020/12/14"
$ws.Range("F22").Value = "Array, Merge"
$ws.Range("G22").Value = "Completed"

$ws.Range("D23").Value = "176. Second Highest Salary"
$ws.Range("D24").Value = "176. Second Highest Salary"

$ws.Range("E21").Value = "2020/12/"
$ws.Range("E24").Value = "2020/12/"

$ws.Range("A23").Value = "LeetCode"
$ws.Range("B23").Value = "Anna"
$ws.Range("C23").Value = "Easy"
$ws.Range("A23").HorizontalAlignment = -4131
$ws.Range("B23").HorizontalAlignment = -4131

$ws.Range("A24").Value = "LeetCode"
$ws.Range("B24").Value = "Stephan"
$ws.Range("C24").Value = "Easy"
$ws.Range("A24").HorizontalAlignment = -4131
$ws.Range("B24").HorizontalAlignment = -4131

# Column D on these rows wraps text (same formatting as the rest of the table);
# column E carries the text-formatted date style; F22 picks up the plain
# applied-font style already used elsewhere in the sheet.
$ws.Range("D21").WrapText = $true
$ws.Range("D22").WrapText = $true
$ws.Range("D23").WrapText = $true
$ws.Range("D24").WrapText = $true
$ws.Range("D25").WrapText = $true

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"

$ws.Rows.Item(21).RowHeight = 14
$ws.Rows.Item(22).RowHeight = 14
$ws.Rows.Item(23).RowHeight = 14
$ws.Rows.Item(24).RowHeight = 14

# Leave the selection where the user's cursor landed after the last edit.
$ws.Range("D24").Select()
